# Added some graphs for the flow meter:
#  - Reposition three of the existing chart graphic frames (Chart 4, Chart 5, Chart 7)
#  - Update the saved cell selection
#  - Re-enter the "time (minutes)" formulas as fill-down ranges so Excel stores
#    them as shared formulas (matching Excel's own behavior when a formula is
#    dragged/filled down a column)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reposition chart graphic frames -----------------------------------
# Chart 4 (2nd chart object): from col6/152400,row33/19050 to col13/457200,row47/95250
$co4 = $ws.ChartObjects(2)
$co4.Left = 362.625
$co4.Top = 496.5
$co4.Width = 433.0625
$co4.Height = 216

# Chart 5 (3rd chart object): from col6/123825,row17/61912 to col13/428625,row31/138112
$co5 = $ws.ChartObjects(3)
$co5.Left = 360.375
$co5.Top = 259.87496062992125
$co5.Width = 433.0625
$co5.Height = 216

# Chart 7 (4th chart object): from col6/138112,row1/138112 to col13/442912,row16/23812
$co7 = $ws.ChartObjects(4)
$co7.Left = 361.49996062992125
$co7.Top = 25.87496062992126
$co7.Width = 433.0625
$co7.Height = 216

# --- Update saved selection ---------------------------------------------
$null = $ws.Range("P35").Select()

# --- Re-apply formulas so they collapse into shared formula groups -----
$ws.Range("E6:E16").Formula = "=B6/60"
$ws.Range("E21:E30").Formula = "=B21/60"
$ws.Range("E34:E44").Formula = "=B34/60"
$ws.Range("E49:E59").Formula = "=B49/60"
